$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jun")

# --- Rename problem titles in column C to "<BOJ number> <name>" ---
$ws.Range("C2").Value  = "10828 스택"
$ws.Range("C3").Value  = "9093 단어 뒤집기"
$ws.Range("C4").Value  = "9012 괄호"
$ws.Range("C5").Value  = "1874 스택 수열"
$ws.Range("C6").Value  = "1406 에디터"
$ws.Range("C7").Value  = "10845 큐"
$ws.Range("C8").Value  = "1158 조세퍼스 문제"
$ws.Range("C9").Value  = "10866 덱"
$ws.Range("C10").Value = "17413 단어 뒤집기 2"
$ws.Range("C11").Value = "10799 쇠막대기"
$ws.Range("C12").Value = "17298 오큰수"
$ws.Range("C13").Value = "17299 오등큰수"
$ws.Range("C14").Value = "1935 후위 표기식 2"
$ws.Range("C15").Value = "1918 후위 표기식"
$ws.Range("C16").Value = "10808 알파벳 개수"
$ws.Range("C17").Value = "10809 알파벳 찾기"
$ws.Range("C18").Value = "10820 문자열 분석"
$ws.Range("C19").Value = "2743 단어 길이 재기"
$ws.Range("C20").Value = "11655 ROT13"
$ws.Range("C21").Value = "10824 네 수"
$ws.Range("C22").Value = "11656 접미사 배열"
$ws.Range("C23").Value = "10430 나머지"
$ws.Range("C24").Value = "2609 최대공약수와 최소공배수"
$ws.Range("C25").Value = "1934 최소공배수"
$ws.Range("C26").Value = "1978 소수 찾기"

# --- Fill in the previously-blank rows 27:29 with new 수학 entries ---
# Copy the formatting of the row above (row 26) down so the new rows keep
# the same borders / fill / date number-format (style index 5 on column D).
$ws.Range("A26:D26").Copy()
$ws.Range("A27:D29").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A27").Value = "수학"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "1929 소수 구하기"
$ws.Range("D27").Value = "12/30/2020"

$ws.Range("A28").Value = "수학"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = "6588 골드바흐의 추측"
$ws.Range("D28").Value = "12/30/2020"

$ws.Range("A29").Value = "수학"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = "10872 팩토리얼"
$ws.Range("D29").Value = "12/30/2020"

# --- Update the saved view state (scroll position / selection) ---
$ws.Activate()
$ws.Range("C29").Select()
